$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -3002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 8000
$ws.Range("J75").Value = 8000
$ws.Range("L75").Value = 8000
$ws.Range("N75").Value = -9748
$ws.Range("H78").Value = 8000
$ws.Range("J78").Value = 8000
$ws.Range("L78").Value = 24000
$ws.Range("N78").Value = -32736

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 10736.429
$ws.Range("I82").Value = 5609.1665
$ws.Range("K82").Value = 5609.1665
$ws.Range("M82").Value = -5226.1665
$ws.Range("H85").Value = 10736.429
$ws.Range("I85").Value = 5609.1665
$ws.Range("K85").Value = 5609.1665
$ws.Range("M85").Value = -4283.1665
$ws.Range("H99").Value = 1253
$ws.Range("I99").Value = 1344.25
$ws.Range("K99").Value = 1344.25
$ws.Range("M99").Value = 153.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 724
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 724
$ws.Range("K2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("M2").Value = 724
$ws.Range("N2").Value = -950
$ws.Range("H3").Value = 1229.8334
$ws.Range("I3").Value = 275
$ws.Range("J3").Value = 1707.25
$ws.Range("K3").Value = 275
$ws.Range("L3").Value = 1707.25
$ws.Range("M3").Value = -162
$ws.Range("N3").Value = -1933.25
$ws.Range("H10").Value = 606.63635
$ws.Range("I10").Value = 176.16667
$ws.Range("J10").Value = 1123.2
$ws.Range("K10").Value = 176.16667
$ws.Range("L10").Value = 1123.2
$ws.Range("M10").Value = -37.16667000000001
$ws.Range("N10").Value = -1401.2
$ws.Range("H11").Value = 721.375
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 721.375
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 721.375
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -1001.375
$ws.Range("H13").Value = 864
$ws.Range("J13").Value = 864
$ws.Range("L13").Value = 864
$ws.Range("N13").Value = -1142
$ws.Range("H31").Value = 3603.7896
$ws.Range("I31").Value = 1864.8
$ws.Range("K31").Value = 1864.8
$ws.Range("M31").Value = -1569.8
$ws.Range("H34").Value = 3603.7896
$ws.Range("I34").Value = 1864.8
$ws.Range("K34").Value = 1864.8
$ws.Range("M34").Value = -1662.8
$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20368
$ws.Range("H58").Value = 5327.375
$ws.Range("I58").Value = 4660
$ws.Range("J58").Value = 9999
$ws.Range("K58").Value = 4660
$ws.Range("L58").Value = 9999
$ws.Range("M58").Value = -4457
$ws.Range("N58").Value = -10405
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490
$ws.Range("H132").Value = 10751.833
$ws.Range("I132").Value = 6503.6665
$ws.Range("K132").Value = 19510.9995
$ws.Range("M132").Value = -16980.9995
$ws.Range("H134").Value = 1812
$ws.Range("I134").Value = 1580.3334
$ws.Range("J134").Value = 2507
$ws.Range("K134").Value = 4741.0002
$ws.Range("L134").Value = 7521
$ws.Range("M134").Value = -2206.0002
$ws.Range("N134").Value = -12591
$ws.Range("H136").Value = 5327.375
$ws.Range("I136").Value = 4660
$ws.Range("J136").Value = 9999
$ws.Range("K136").Value = 13980
$ws.Range("L136").Value = 29997
$ws.Range("M136").Value = -11430
$ws.Range("N136").Value = -35097

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 863.3333
$ws.Range("I26").Value = 863.3333
$ws.Range("K26").Value = 2589.9999
$ws.Range("M26").Value = -2301.9999
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30540
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31872
$ws.Range("H97").Value = 77
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H137").Value = 1800
$ws.Range("I137").Value = 1800
$ws.Range("K137").Value = 5400
$ws.Range("M137").Value = -300

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1498
$ws.Range("J4").Value = 1498
$ws.Range("L4").Value = 1498
$ws.Range("N4").Value = -1722
$ws.Range("H10").Value = 4301
$ws.Range("J10").Value = 1450
$ws.Range("L10").Value = 1450
$ws.Range("N10").Value = -1788
$ws.Range("H11").Value = 7166983.5
$ws.Range("I11").Value = 8250000
$ws.Range("K11").Value = 8250000
$ws.Range("M11").Value = -8249861
$ws.Range("H13").Value = 404.85715
$ws.Range("I13").Value = 162.5
$ws.Range("J13").Value = 501.8
$ws.Range("K13").Value = 162.5
$ws.Range("L13").Value = 501.8
$ws.Range("M13").Value = -23.5
$ws.Range("N13").Value = -779.8
$ws.Range("H17").Value = 1029
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 2000
$ws.Range("N17").Value = -2336
$ws.Range("H22").Value = 4000
$ws.Range("J22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("N22").Value = -5058
$ws.Range("H25").Value = 5949.75
$ws.Range("J25").Value = 5949.75
$ws.Range("L25").Value = 5949.75
$ws.Range("N25").Value = -7007.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H3").Value = 4000
$ws.Range("J3").Value = 4000
$ws.Range("L3").Value = 4000
$ws.Range("N3").Value = -4224
$ws.Range("H10").Value = 2590.8
$ws.Range("J10").Value = 3634.6667
$ws.Range("L10").Value = 3634.6667
$ws.Range("N10").Value = -3914.6667
$ws.Range("H12").Value = 3111.111
$ws.Range("J12").Value = 1133.3334
$ws.Range("L12").Value = 1133.3334
$ws.Range("N12").Value = -1473.3334
$ws.Range("H13").Value = 18999.5
$ws.Range("I13").Value = 18999.5
$ws.Range("K13").Value = 18999.5
$ws.Range("M13").Value = -18859.5
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2344
$ws.Range("H15").Value = 4000
$ws.Range("J15").Value = 4000
$ws.Range("L15").Value = 4000
$ws.Range("N15").Value = -4340
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("N17").Value = 0
$ws.Range("H19").Value = 500
$ws.Range("J19").Value = 500
$ws.Range("L19").Value = 500
$ws.Range("N19").Value = -840
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0
$ws.Range("H136").Value = 9998
$ws.Range("I136").Value = 10004
$ws.Range("K136").Value = 30012
$ws.Range("M136").Value = -27462
